$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear Priority column (P0) for existing rows 4-7 ---
$ws.Range("E4:E7").ClearContents()

# --- 2. Row 7 gains a "Number of Test Cases" value (F7 = 7) ---
$ws.Range("F6").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F7").Value = 7

# --- 3. Update row 6 text (Home Page Currency -> Currency) ---
$ws.Range("B6").Value = "(TS_003)" + [char]10 + "Currency"
$ws.Range("D6").Value = "Validate the working of home page > Currency"

# --- 4. Update row 7 text (Home Page Contact Us -> Contact Us) ---
$ws.Range("B7").Value = "(TS_004)" + [char]10 + "Contact Us"
$ws.Range("D7").Value = "Validate the working of home page > Contact Us"

# --- 5. Build new row 8 (B, C, D, F only - no E) ---
$ws.Range("B6").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("C6").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("F6").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("B8,C8,D8,F8").ClearContents()
$ws.Range("B8").Value = "(TS_005)" + [char]10 + "My Account"
$ws.Range("C8").Value = "FRS"
$ws.Range("D8").Value = "Validate the working of home page > My Account"
$ws.Range("F8").Value = 4
$ws.Rows.Item(8).RowHeight = 30

# --- 6. Build new row 9 (B, C, D only) ---
$ws.Range("B6").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C6").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("B9,C9,D9").ClearContents()
$ws.Range("B9").Value = "(TS_006)" + [char]10 + " Account"
$ws.Range("C9").Value = "FRS"
$ws.Range("D9").Value = "Validate the working of home page > My Account> Account"
$ws.Rows.Item(9).RowHeight = 30

# --- 7. Add rows 10 and 11 with just a formatted D cell ---
$ws.Range("D6").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D10").ClearContents()
$ws.Range("D6").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D11").ClearContents()
$ws.Rows.Item(10).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 15

# --- 8. Adjust column D width (72 displayed) ---
$ws.Columns.Item(4).ColumnWidth = 71.16666666666667

# --- 9. Update selection ---
$ws.Range("D10").Select()
